# Scheduled runner: refresh Leve profit calculations (currentAveragePrice* / LevePrice* / LeveProfit*)
# Updates reflect latest market-board pricing snapshot per Leve row, across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 315.11765
$ws.Cells.Item(2, 9).Value = 292.30768
$ws.Cells.Item(2, 10).Value = 389.25
$ws.Cells.Item(2, 11).Value = 292.30768
$ws.Cells.Item(2, 12).Value = 389.25
$ws.Cells.Item(2, 13).Value = -179.30768
$ws.Cells.Item(2, 14).Value = -615.25

$ws.Cells.Item(17, 8).Value = 9489.846
$ws.Cells.Item(17, 10).Value = 10155.667
$ws.Cells.Item(17, 12).Value = 30467.001
$ws.Cells.Item(17, 14).Value = -30803.001

$ws.Cells.Item(28, 8).Value = 1440.4166
$ws.Cells.Item(28, 9).Value = 1051
$ws.Cells.Item(28, 10).Value = 2386.1428
$ws.Cells.Item(28, 11).Value = 1051
$ws.Cells.Item(28, 12).Value = 2386.1428
$ws.Cells.Item(28, 13).Value = -566
$ws.Cells.Item(28, 14).Value = -3356.1428

$ws.Cells.Item(41, 8).Value = 117.5
$ws.Cells.Item(41, 9).Value = 86.25
$ws.Cells.Item(41, 10).Value = 180
$ws.Cells.Item(41, 11).Value = 86.25
$ws.Cells.Item(41, 12).Value = 180
$ws.Cells.Item(41, 13).Value = 353.75
$ws.Cells.Item(41, 14).Value = -1060

$ws.Cells.Item(55, 8).Value = 467.14706
$ws.Cells.Item(55, 9).Value = 403.375
$ws.Cells.Item(55, 10).Value = 523.8333
$ws.Cells.Item(55, 11).Value = 403.375
$ws.Cells.Item(55, 12).Value = 523.8333
$ws.Cells.Item(55, 13).Value = -189.375
$ws.Cells.Item(55, 14).Value = -951.8333

$ws.Cells.Item(86, 8).Value = 5245.5
$ws.Cells.Item(86, 10).Value = 5535.4287
$ws.Cells.Item(86, 12).Value = 5535.4287
$ws.Cells.Item(86, 14).Value = -7781.4287

$ws.Cells.Item(89, 8).Value = 5245.5
$ws.Cells.Item(89, 10).Value = 5535.4287
$ws.Cells.Item(89, 12).Value = 27677.1435
$ws.Cells.Item(89, 14).Value = -38909.14350000001

$ws.Cells.Item(100, 8).Value = 2763.1875
$ws.Cells.Item(100, 9).Value = 1415.3334
$ws.Cells.Item(100, 10).Value = 4496.143
$ws.Cells.Item(100, 11).Value = 1415.3334
$ws.Cells.Item(100, 12).Value = 4496.143
$ws.Cells.Item(100, 13).Value = -874.3334
$ws.Cells.Item(100, 14).Value = -5578.143

$ws.Cells.Item(132, 8).Value = 1465.174
$ws.Cells.Item(132, 9).Value = 1250.05
$ws.Cells.Item(132, 10).Value = 2899.3333
$ws.Cells.Item(132, 11).Value = 3750.15
$ws.Cells.Item(132, 12).Value = 8697.999899999999
$ws.Cells.Item(132, 13).Value = -1220.15
$ws.Cells.Item(132, 14).Value = -13757.9999

$ws.Cells.Item(137, 8).Value = 3194.8306
$ws.Cells.Item(137, 9).Value = 2536.75
$ws.Cells.Item(137, 10).Value = 3362.851
$ws.Cells.Item(137, 11).Value = 7610.25
$ws.Cells.Item(137, 12).Value = 10088.553
$ws.Cells.Item(137, 13).Value = -5060.25
$ws.Cells.Item(137, 14).Value = -15188.553

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5643.5312
$ws.Cells.Item(32, 9).Value = 4365.6553
$ws.Cells.Item(32, 11).Value = 4365.6553
$ws.Cells.Item(32, 13).Value = -4078.6553

$ws.Cells.Item(45, 8).Value = 3080.9167
$ws.Cells.Item(45, 9).Value = 2947.1
$ws.Cells.Item(45, 10).Value = 3750
$ws.Cells.Item(45, 11).Value = 2947.1
$ws.Cells.Item(45, 12).Value = 3750
$ws.Cells.Item(45, 13).Value = -2570.1
$ws.Cells.Item(45, 14).Value = -4504

$ws.Cells.Item(61, 8).Value = 5390.24
$ws.Cells.Item(61, 9).Value = 4405.9165
$ws.Cells.Item(61, 11).Value = 4405.9165
$ws.Cells.Item(61, 13).Value = -4193.9165

$ws.Cells.Item(74, 8).Value = 37041564
$ws.Cells.Item(74, 9).Value = 41669384
$ws.Cells.Item(74, 11).Value = 41669384
$ws.Cells.Item(74, 13).Value = -41668510

$ws.Cells.Item(77, 8).Value = 37041564
$ws.Cells.Item(77, 9).Value = 41669384
$ws.Cells.Item(77, 11).Value = 208346920
$ws.Cells.Item(77, 13).Value = -208342552

$ws.Cells.Item(119, 8).Value = 69317
$ws.Cells.Item(119, 10).Value = 69317
$ws.Cells.Item(119, 12).Value = 69317
$ws.Cells.Item(119, 14).Value = -78993

$ws.Cells.Item(122, 8).Value = 3639.913
$ws.Cells.Item(122, 9).Value = 3165.7646
$ws.Cells.Item(122, 11).Value = 9497.293799999999
$ws.Cells.Item(122, 13).Value = -7047.293799999999

$ws.Cells.Item(132, 8).Value = 3176.5833
$ws.Cells.Item(132, 9).Value = 2249
$ws.Cells.Item(132, 11).Value = 6747
$ws.Cells.Item(132, 13).Value = -4217

$ws.Cells.Item(136, 8).Value = 5390.24
$ws.Cells.Item(136, 9).Value = 4405.9165
$ws.Cells.Item(136, 11).Value = 13217.7495
$ws.Cells.Item(136, 13).Value = -10667.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 26773.273
$ws.Cells.Item(105, 9).Value = 26563.375
$ws.Cells.Item(105, 10).Value = 27333
$ws.Cells.Item(105, 11).Value = 26563.375
$ws.Cells.Item(105, 12).Value = 27333
$ws.Cells.Item(105, 13).Value = -24816.375
$ws.Cells.Item(105, 14).Value = -30827

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 50911.137
$ws.Cells.Item(31, 9).Value = 1751.4166
$ws.Cells.Item(31, 11).Value = 1751.4166
$ws.Cells.Item(31, 13).Value = -1456.4166

$ws.Cells.Item(34, 8).Value = 50911.137
$ws.Cells.Item(34, 9).Value = 1751.4166
$ws.Cells.Item(34, 11).Value = 1751.4166
$ws.Cells.Item(34, 13).Value = -1549.4166

$ws.Cells.Item(58, 8).Value = 11432.538
$ws.Cells.Item(58, 9).Value = 2405.5
$ws.Cells.Item(58, 10).Value = 15444.556
$ws.Cells.Item(58, 11).Value = 2405.5
$ws.Cells.Item(58, 12).Value = 15444.556
$ws.Cells.Item(58, 13).Value = -2202.5
$ws.Cells.Item(58, 14).Value = -15850.556

$ws.Cells.Item(99, 8).Value = 5033
$ws.Cells.Item(99, 9).Value = 4999.5
$ws.Cells.Item(99, 10).Value = 5100
$ws.Cells.Item(99, 11).Value = 4999.5
$ws.Cells.Item(99, 12).Value = 5100
$ws.Cells.Item(99, 13).Value = -3501.5
$ws.Cells.Item(99, 14).Value = -8096

$ws.Cells.Item(126, 8).Value = 5033
$ws.Cells.Item(126, 9).Value = 4999.5
$ws.Cells.Item(126, 10).Value = 5100
$ws.Cells.Item(126, 11).Value = 14998.5
$ws.Cells.Item(126, 12).Value = 15300
$ws.Cells.Item(126, 13).Value = -12528.5
$ws.Cells.Item(126, 14).Value = -20240

$ws.Cells.Item(132, 8).Value = 2765.0605
$ws.Cells.Item(132, 9).Value = 1927.6538
$ws.Cells.Item(132, 10).Value = 5875.4287
$ws.Cells.Item(132, 11).Value = 5782.9614
$ws.Cells.Item(132, 12).Value = 17626.2861
$ws.Cells.Item(132, 13).Value = -3252.9614
$ws.Cells.Item(132, 14).Value = -22686.2861

$ws.Cells.Item(134, 8).Value = 3009.4666
$ws.Cells.Item(134, 9).Value = 2216.6924
$ws.Cells.Item(134, 11).Value = 6650.0772
$ws.Cells.Item(134, 13).Value = -4115.0772

$ws.Cells.Item(136, 8).Value = 11432.538
$ws.Cells.Item(136, 9).Value = 2405.5
$ws.Cells.Item(136, 10).Value = 15444.556
$ws.Cells.Item(136, 11).Value = 7216.5
$ws.Cells.Item(136, 12).Value = 46333.66800000001
$ws.Cells.Item(136, 13).Value = -4666.5
$ws.Cells.Item(136, 14).Value = -51433.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 1148
$ws.Cells.Item(26, 9).Value = 1645
$ws.Cells.Item(26, 10).Value = 849.8
$ws.Cells.Item(26, 11).Value = 4935
$ws.Cells.Item(26, 12).Value = 2549.4
$ws.Cells.Item(26, 13).Value = -4647
$ws.Cells.Item(26, 14).Value = -3125.4

$ws.Cells.Item(59, 8).Value = 1773
$ws.Cells.Item(59, 9).Value = 1958.3334
$ws.Cells.Item(59, 11).Value = 5875.0002
$ws.Cells.Item(59, 13).Value = -5335.0002

$ws.Cells.Item(61, 8).Value = 1015.9167
$ws.Cells.Item(61, 9).Value = 126.5
$ws.Cells.Item(61, 10).Value = 1905.3334
$ws.Cells.Item(61, 11).Value = 379.5
$ws.Cells.Item(61, 12).Value = 5716.0002
$ws.Cells.Item(61, 13).Value = -164.5
$ws.Cells.Item(61, 14).Value = -6146.0002

$ws.Cells.Item(131, 8).Value = 9128838
$ws.Cells.Item(131, 10).Value = 6274358.5
$ws.Cells.Item(131, 12).Value = 18823075.5
$ws.Cells.Item(131, 14).Value = -18833155.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1639.1923
$ws.Cells.Item(97, 9).Value = 1209.4286
$ws.Cells.Item(97, 11).Value = 1209.4286
$ws.Cells.Item(97, 13).Value = -713.4286

$ws.Cells.Item(132, 8).Value = 2552.682
$ws.Cells.Item(132, 9).Value = 1366.5
$ws.Cells.Item(132, 10).Value = 4628.5
$ws.Cells.Item(132, 11).Value = 4099.5
$ws.Cells.Item(132, 12).Value = 13885.5
$ws.Cells.Item(132, 13).Value = -1569.5
$ws.Cells.Item(132, 14).Value = -18945.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1687.8
$ws.Cells.Item(16, 9).Value = 1062.9231
$ws.Cells.Item(16, 11).Value = 1062.9231
$ws.Cells.Item(16, 13).Value = -892.9231

$ws.Cells.Item(44, 8).Value = 100000
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 13).ClearContents()

$ws.Cells.Item(46, 8).Value = 5043.1035
$ws.Cells.Item(46, 9).Value = 3521.3572
$ws.Cells.Item(46, 10).Value = 6463.4
$ws.Cells.Item(46, 11).Value = 3521.3572
$ws.Cells.Item(46, 12).Value = 6463.4
$ws.Cells.Item(46, 13).Value = -3333.3572
$ws.Cells.Item(46, 14).Value = -6839.4

$ws.Cells.Item(74, 8).Value = 62499.5
$ws.Cells.Item(74, 10).Value = 74999
$ws.Cells.Item(74, 12).Value = 74999
$ws.Cells.Item(74, 14).Value = -76995

$ws.Cells.Item(77, 8).Value = 62499.5
$ws.Cells.Item(77, 10).Value = 74999
$ws.Cells.Item(77, 12).Value = 224997
$ws.Cells.Item(77, 14).Value = -234981

$ws.Cells.Item(92, 8).Value = 47713.855
$ws.Cells.Item(92, 10).Value = 47713.855
$ws.Cells.Item(92, 12).Value = 47713.855
$ws.Cells.Item(92, 14).Value = -52705.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 19632.666
$ws.Cells.Item(58, 9).Value = 16699
$ws.Cells.Item(58, 11).Value = 16699
$ws.Cells.Item(58, 13).Value = -16391

$ws.Cells.Item(119, 8).Value = 59822.5
$ws.Cells.Item(119, 10).Value = 59822.5
$ws.Cells.Item(119, 12).Value = 59822.5
$ws.Cells.Item(119, 14).Value = -69498.5

$ws.Cells.Item(132, 8).Value = 2121
$ws.Cells.Item(132, 9).Value = 1182.8
$ws.Cells.Item(132, 10).Value = 11503
$ws.Cells.Item(132, 11).Value = 3548.4
$ws.Cells.Item(132, 12).Value = 34509
$ws.Cells.Item(132, 13).Value = -1018.4
$ws.Cells.Item(132, 14).Value = -39569

$ws.Cells.Item(136, 8).Value = 3607.739
$ws.Cells.Item(136, 9).Value = 2400
$ws.Cells.Item(136, 10).Value = 5872.25
$ws.Cells.Item(136, 11).Value = 7200
$ws.Cells.Item(136, 12).Value = 17616.75
$ws.Cells.Item(136, 13).Value = -4650
$ws.Cells.Item(136, 14).Value = -22716.75

